$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: update the teacher's details (previously "Ahmad Zidan" / "zidan" / "Zidan456" / "ADMIN")
$ws.Range("A4").Value = 2832083
$ws.Range("B4").Value = "guru matematika"
$ws.Range("C4").Value = "bufit"
$ws.Range("D4").Value = "ireireir"
$ws.Range("E4").Value = "Guru"

# Row 5: new record
$ws.Range("A5").Value = 12345
$ws.Range("B5").Value = "Alvin Ganteng"
$ws.Range("C5").Value = "Alvin"
$ws.Range("D5").Value = 12345
$ws.Range("E5").Value = "XI-RPL-2"

# Row 6: new record
$ws.Range("A6").Value = 54321
$ws.Range("B6").Value = "Anton"
$ws.Range("C6").Value = "Anton"
$ws.Range("D6").Value = 12345
$ws.Range("E6").Value = "XI-RPL-2"

# Make sure all data cells carry the same style as the rest of the table (s="1")
$ws.Range("A2:E2").Copy()
$ws.Range("A4:E6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
